# Add a new "Llamma3 Fine-tunded" column of results to the comparison table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H; this shifts the existing H:P data
# (and their formatting) one column to the right, into I:Q.
$ws.Columns("H").Insert()

# Match the width of the newly inserted column to the rest of the table
# (closest width achievable through the character-width column API).
$ws.Columns("H").ColumnWidth = 23.142857142857142

# New header for column H.
$ws.Range("H1").Value = "Llamma3 Fine-tunded"

# New per-model-row values for column H.
$ws.Range("H2").Value = 0.368
$ws.Range("H3").Value = 0.4
$ws.Range("H4").Value = 0.45
$ws.Range("H5").Value = 0.2
$ws.Range("H6").Value = 0.6
$ws.Range("H7").Value = 0.4
$ws.Range("H8").Value = 0.4

# Row 3's new cell needs the two-decimal percent format (style used by the
# rest of that row), which differs from the format inherited from column G.
$ws.Range("H3").NumberFormat = "0.00%"

# Update the active cell/selection to match the author's final cursor position.
$null = $ws.Range("J15").Select()
